$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Publication date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was previously blank)
$meta.Range("B9").Value = "Alvearie Team"

# "Contact / No display for ContactDetail" row becomes "Jurisdiction / United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The old sheet had a duplicate "Contact" row right below (row 11) - remove it,
# shifting every following row up by one.
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements" ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition now describe this specific extension
$elements.Range("K2").Value = "Documented Sex"
$elements.Range("L2").Value = "Sex coding indicated on the legal document"
